# Regenerate save_data column G ("K" - strikeouts) to replace the old
# "Strike#" values with the recalculated K values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 2
    6  = 1
    7  = 0
    9  = 1
    10 = 2
    11 = 3
    14 = 1
    17 = 1
    18 = 1
    19 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
